# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# on the crypto listing sheet to reflect refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the target cells keep a plain "Text" format so the values
# are stored as strings (matching the original inlineStr cells),
# e.g. "39.90" does not get silently rewritten to "39.9" and
# "6.30%" keeps its literal percent-sign text rather than becoming 0.063.
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '6.30%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '39.90'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '7.88%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.732'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '11.70%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08104'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '4.04%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.570'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '3.85%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '8.685'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '4.54%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.966'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '5.12%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9425'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.07%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1283'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '15.68%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1983'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '5.66%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09214'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '4.59%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03463'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '4.94%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09615'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.40%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001337'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-3.14%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006018'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-3.06%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.58%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3528'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.21%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.576'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '18.54%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1411'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '9.17%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2428'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2.23%'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.18%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001252'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '4.14%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004326'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '1.20%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-15.09%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003992'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '37.39%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02522'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '17.70%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05219'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '5.42%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007314'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-3.64%'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '5.56%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.008997'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '5.56%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002191'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '10.04%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01001'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '16.30%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006694'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '1.60%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.12%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-12.85%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.001801'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '24.58%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.12%'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.12%'

Write-Output "Updated $([string]68) cells in columns D and E"
